# Update Work Week and Social Spending
# (commit message is generic/boilerplate; the actual change updates the
# "GDP per Capita" data series on the "Data" sheet: revises the figures for
# 1973 and 1980-2010, and appends newly available years 2011-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The "Data" column (E) stores its numbers as *text* shared-strings in the
# source workbook (not native numeric cells). Assigning a numeric-looking
# string straight to .Value lets Excel auto-coerce it to a number (and can
# round high-precision decimals), so we briefly force Text number format
# before writing, then clear the formatting back off again afterwards so we
# don't leave a lingering custom style on the cells.
# (Multi-area Range("A,B:C") unions aren't reliable here, so the two
# contiguous blocks are formatted separately.)
$textCellsA = $ws.Range("E2")
$textCellsB = $ws.Range("E9:E45")
$textCellsA.NumberFormat = "@"
$textCellsB.NumberFormat = "@"

# --- Revised figures for existing years ---
$ws.Range("E2").Value = "9456"
$ws.Range("E9").Value = "10519"
$ws.Range("E10").Value = "11053"
$ws.Range("E11").Value = "11611"
$ws.Range("E12").Value = "12234"
$ws.Range("E13").Value = "13128"
$ws.Range("E14").Value = "13868"
$ws.Range("E15").Value = "13839"
$ws.Range("E16").Value = "13707"
$ws.Range("E17").Value = "14665"
$ws.Range("E18").Value = "14161"
$ws.Range("E19").Value = "12140"
$ws.Range("E20").Value = "9535.05058084973"
$ws.Range("E21").Value = "5286.92481701539"
$ws.Range("E22").Value = "3792.61767813275"
$ws.Range("E23").Value = "3463.47608288493"
$ws.Range("E24").Value = "3610.84075821328"
$ws.Range("E25").Value = "4030.88466295627"
$ws.Range("E26").Value = "4480.89785963119"
$ws.Range("E27").Value = "4619.90479030964"
$ws.Range("E28").Value = "4763.3505565953"
$ws.Range("E29").Value = "4891.66405370768"
$ws.Range("E30").Value = "5125.12737032545"
$ws.Range("E31").Value = "5406.95332679979"
$ws.Range("E32").Value = "6013.00915262526"
$ws.Range("E33").Value = "6342.70816986243"
$ws.Range("E34").Value = "6849.78904677909"
$ws.Range("E35").Value = "7398.05264239695"
$ws.Range("E36").Value = "8300.14388042897"
$ws.Range("E37").Value = "8461.91734250787"
$ws.Range("E38").Value = "8062.08378622631"
$ws.Range("E39").Value = "8443.4360888918"

# --- Newly published years 2011-2016 ---
$ws.Range("A40").Value = 268
$ws.Range("B40").Value = "Georgia"
$ws.Range("C40").Value = "GDP per Capita"
$ws.Range("D40").Value = 2011
$ws.Range("E40").Value = "8946"

$ws.Range("A41").Value = 268
$ws.Range("B41").Value = "Georgia"
$ws.Range("C41").Value = "GDP per Capita"
$ws.Range("D41").Value = 2012
$ws.Range("E41").Value = "9491"

$ws.Range("A42").Value = 268
$ws.Range("B42").Value = "Georgia"
$ws.Range("C42").Value = "GDP per Capita"
$ws.Range("D42").Value = 2013
$ws.Range("E42").Value = "9829"

$ws.Range("A43").Value = 268
$ws.Range("B43").Value = "Georgia"
$ws.Range("C43").Value = "GDP per Capita"
$ws.Range("D43").Value = 2014
$ws.Range("E43").Value = "10297"

$ws.Range("A44").Value = 268
$ws.Range("B44").Value = "Georgia"
$ws.Range("C44").Value = "GDP per Capita"
$ws.Range("D44").Value = 2015
$ws.Range("E44").Value = "10603"

$ws.Range("A45").Value = 268
$ws.Range("B45").Value = "Georgia"
$ws.Range("C45").Value = "GDP per Capita"
$ws.Range("D45").Value = 2016
$ws.Range("E45").Value = "10899"

# Drop the temporary Text format so the cells end up back on the default
# (General) style, same as the rest of the sheet.
$textCellsA.ClearFormats()
$textCellsB.ClearFormats()
